$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp text
$ws.Range("A1").Value = "Datos actualizados a 4 de Julio de 2020 a las 04:48"

# 2. Alemania (row 18): update Recuperados / Casos activos
$ws.Range("D18").Value = 181300
$ws.Range("E18").Value = 6627

# 3. Swap Bolivia / Panama (rows 43 & 44) with updated figures
$ws.Range("A43").Value = "Bolivia"
$ws.Range("B43").Value = 36818
$ws.Range("C43").Value = 1290
$ws.Range("D43").Value = 10766
$ws.Range("E43").Value = 24732
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 49
$ws.Range("H43").Value = 1320

$ws.Range("A44").Value = "Panama"
$ws.Range("B44").Value = 35995
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 16945
$ws.Range("E44").Value = 18352
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 698

# 4. Swap Guatemala / Honduras (rows 55 & 56) with updated figures
$ws.Range("A55").Value = "Honduras"
$ws.Range("B55").Value = 22116
$ws.Range("C55").Value = 996
$ws.Range("D55").Value = 2250
$ws.Range("E55").Value = 19261
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 14
$ws.Range("H55").Value = 605

$ws.Range("A56").Value = "Guatemala"
$ws.Range("B56").Value = 21293
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 3315
$ws.Range("E56").Value = 17098
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 880

# 5. Australia (row 75): update figures
$ws.Range("B75").Value = 8260
$ws.Range("C75").Value = 5
$ws.Range("D75").Value = 7319
$ws.Range("E75").Value = 837

# 6. Swap San Vicente y las Granadinas / Belice (rows 197 & 198) with updated figures
$ws.Range("A197").Value = "Belice"
$ws.Range("B197").Value = 30
$ws.Range("C197").Value = 2
$ws.Range("D197").Value = 19
$ws.Range("E197").Value = 9
$ws.Range("F197").Value = 0
$ws.Range("G197").Value = 0
$ws.Range("H197").Value = 2

$ws.Range("A198").Value = "San Vicente y las Granadinas"
$ws.Range("B198").Value = 29
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 29
$ws.Range("E198").Value = 0
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 0
